$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'25.858.62"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "'  +0.42%  "
$ws.Range("E2").ClearFormats()
$ws.Range("D3").Value = "'1.635.15"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "'  +0.03%  "
$ws.Range("E3").ClearFormats()
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "'  -0.28%  "
$ws.Range("E4").ClearFormats()
$ws.Range("D5").Value = "'214.69"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "'  -0.36%  "
$ws.Range("E5").ClearFormats()
$ws.Range("E6").Value = "'  -0.45%  "
$ws.Range("E6").ClearFormats()
$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "'  -0.31%  "
$ws.Range("E7").ClearFormats()
$ws.Range("E8").Value = "'  -0.87%  "
$ws.Range("E8").ClearFormats()
$ws.Range("E9").Value = "'  -0.71%  "
$ws.Range("E9").ClearFormats()
$ws.Range("D10").Value = "'19.63"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "'  +0.23%  "
$ws.Range("E10").ClearFormats()
$ws.Range("E11").Value = "'  +0.42%  "
$ws.Range("E11").ClearFormats()
$ws.Range("E12").Value = "'  +0.68%  "
$ws.Range("E12").ClearFormats()
$ws.Range("D13").Value = "'1.860.04"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "'  -0.01%  "
$ws.Range("E13").ClearFormats()
$ws.Range("D14").Value = "'1.618.67"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "'  -1.20%  "
$ws.Range("E14").ClearFormats()
$ws.Range("D15").Value = "'0.551"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "'  -1.01%  "
$ws.Range("E15").ClearFormats()
$ws.Range("E16").Value = "'  -0.75%  "
$ws.Range("E16").ClearFormats()
$ws.Range("D17").Value = "'62.88"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "'  +0.22%  "
$ws.Range("E17").ClearFormats()
$ws.Range("D18").Value = "'25.832.44"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "'  +0.25%  "
$ws.Range("E18").ClearFormats()
$ws.Range("D19").Value = "'1.00"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "'  -0.25%  "
$ws.Range("E19").ClearFormats()
$ws.Range("D20").Value = "'4.44"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "'  -0.13%  "
$ws.Range("E20").ClearFormats()
$ws.Range("D21").Value = "'191.33"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "'  -1.30%  "
$ws.Range("E21").ClearFormats()
$ws.Range("D22").Value = "'9.97"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "'  +0.23%  "
$ws.Range("E22").ClearFormats()
$ws.Range("E23").Value = "'  +0.69%  "
$ws.Range("E23").ClearFormats()
$ws.Range("D24").Value = "'0.998"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "'  -0.38%  "
$ws.Range("E24").ClearFormats()
$ws.Range("D25").Value = "'1.81"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "'  -1.91%  "
$ws.Range("E25").ClearFormats()
$ws.Range("D26").Value = "'142.40"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "'  +0.87%  "
$ws.Range("E26").ClearFormats()
$ws.Range("E27").Value = "'  +0.79%  "
$ws.Range("E27").ClearFormats()
$ws.Range("D28").Value = "'6.85"
$ws.Range("D28").ClearFormats()
$ws.Range("D29").Value = "'15.51"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "'  -0.09%  "
$ws.Range("E29").ClearFormats()
$ws.Range("E30").Value = "'  -0.46%  "
$ws.Range("E30").ClearFormats()
$ws.Range("D31").Value = "'0.0493"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "'  +0.18%  "
$ws.Range("E31").ClearFormats()
$ws.Range("D32").Value = "'3.33"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "'  +0.05%  "
$ws.Range("E32").ClearFormats()
$ws.Range("D33").Value = "'3.24"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "'  -0.42%  "
$ws.Range("E33").ClearFormats()
$ws.Range("E34").Value = "'  +0.65%  "
$ws.Range("E34").ClearFormats()
$ws.Range("D35").Value = "'2.40"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "'  +0.49%  "
$ws.Range("E35").ClearFormats()
$ws.Range("D36").Value = "'0.906"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "'  +0.52%  "
$ws.Range("E36").ClearFormats()
$ws.Range("D37").Value = "'1.147.19"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "'  +2.14%  "
$ws.Range("E37").ClearFormats()
$ws.Range("D38").Value = "'0.544"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "'  -0.80%  "
$ws.Range("E38").ClearFormats()
$ws.Range("E39").Value = "'  -1.05%  "
$ws.Range("E39").ClearFormats()
$ws.Range("E40").Value = "'  +0.48%  "
$ws.Range("E40").ClearFormats()
$ws.Range("D41").Value = "'0.999"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "'  -0.29%  "
$ws.Range("E41").ClearFormats()
$ws.Range("D42").Value = "'5.62"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "'  +0.85%  "
$ws.Range("E42").ClearFormats()
$ws.Range("D43").Value = "'100.53"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "'  +0.84%  "
$ws.Range("E43").ClearFormats()
$ws.Range("D44").Value = "'0.804"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "'  +0.11%  "
$ws.Range("E44").ClearFormats()
$ws.Range("D45").Value = "'1.770.30"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "'  +0.04%  "
$ws.Range("E45").ClearFormats()
$ws.Range("B46").Value = "'BabyDogeCoin"
$ws.Range("B46").ClearFormats()
$ws.Range("C46").Value = "'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("C46").ClearFormats()
$ws.Range("D46").Value = "'0.0₆0111"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "'  -1.21%  "
$ws.Range("E46").ClearFormats()
$ws.Range("B47").Value = "'Aave"
$ws.Range("B47").ClearFormats()
$ws.Range("C47").Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("C47").ClearFormats()
$ws.Range("D47").Value = "'55.65"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "'  +1.03%  "
$ws.Range("E47").ClearFormats()
$ws.Range("B48").Value = "'Cronos"
$ws.Range("B48").ClearFormats()
$ws.Range("C48").Value = "'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("C48").ClearFormats()
$ws.Range("D48").Value = "'0.0512"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "'  +2.24%  "
$ws.Range("E48").ClearFormats()
$ws.Range("B49").Value = "'RenderToken"
$ws.Range("B49").ClearFormats()
$ws.Range("C49").Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("C49").ClearFormats()
$ws.Range("D49").Value = "'1.46"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "'  +5.35%  "
$ws.Range("E49").ClearFormats()
$ws.Range("B50").Value = "'Mantle"
$ws.Range("B50").ClearFormats()
$ws.Range("C50").Value = "'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("C50").ClearFormats()
$ws.Range("D50").Value = "'0.417"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "'  -0.02%  "
$ws.Range("E50").ClearFormats()
$ws.Range("B51").Value = "'EnergySwap"
$ws.Range("B51").ClearFormats()
$ws.Range("C51").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("C51").ClearFormats()
$ws.Range("D51").Value = "'7.55"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "'  -0.15%  "
$ws.Range("E51").ClearFormats()

Write-Host "Update complete"
